$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K")

$rowData = @{}
$rowData[2] = @(-17.39050486507988, 1.974879924760914, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[3] = @(-17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, 3.123389435717151, -17.39050486507988, -17.39050486507988)
$rowData[4] = @(-17.39050486507988, 2.219245610660054, 2.23563371211787, -17.39050486507988, 3.412040687108591, -17.39050486507988, 1.740440177932122, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[5] = @(-17.39050486507988, 1.959374568989561, -17.39050486507988, -17.39050486507988, -17.39050486507988, 2.923875910819471, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[6] = @(-17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[7] = @(2.622104671092812, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[8] = @(-17.39050486507988, -17.39050486507988, -17.39050486507988, 1.751854023462338, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[9] = @(3.791128467178512, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[10] = @(-17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, 1.461893579591917, -17.39050486507988, 1.969058118911824)
$rowData[11] = @(-17.39050486507988, -17.39050486507988, -17.39050486507988, 2.888798517549762, -17.39050486507988, 2.545386748573573, -17.39050486507988, -17.39050486507988, -17.39050486507988, 1.561294924837557)
$rowData[12] = @(-17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[13] = @(-17.39050486507988, -17.39050486507988, -17.39050486507988, 2.431208401366522, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, 4.321920117958101, 1.583598784881665)
$rowData[14] = @(-17.39050486507988, -17.39050486507988, 1.331372003742765, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, 2.075804648088508)
$rowData[15] = @(-17.39050486507988, -17.39050486507988, 1.119088461283763, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[16] = @(-17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988)
$rowData[17] = @(-17.39050486507988, 1.808565943907243, 2.16758881599518, -17.39050486507988, -17.39050486507988, -17.39050486507988, 0.6584981074933609, 0.7074500872490243, -17.39050486507988, -17.39050486507988)
$rowData[18] = @(-17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, -17.39050486507988, 1.201795215669559, 0.7537998693806095, -17.39050486507988, -17.39050486507988)
$rowData[19] = @(-17.39050486507988, -17.39050486507988, 1.723606004522762, -17.39050486507988, -17.39050486507988, -17.39050486507988, 1.78312591765526, 1.750457177289984, -17.39050486507988, -17.39050486507988)
$rowData[20] = @(-17.39050486507988, 0.8182676976824615, 1.488710413043474, -17.39050486507988, 3.22579358879141, -17.39050486507988, 2.251520378008469, 0.8852604578714924, -17.39050486507988, 2.565753585107771)
$rowData[21] = @(-17.39050486507988, 1.174600933095772, -17.39050486507988, 1.938102136328417, -17.39050486507988, 2.71664706581709, 2.194180602321346, -17.39050486507988, -17.39050486507988, -17.39050486507988)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$r"
        $ws.Range($addr).Value = $vals[$i]
    }
}

Write-Host "applied PSSM updates for rows 2-21"
